# Update "F" column (想去人数 / want-to-go count) values on three sheets:
# 展览 (Exhibitions), 演出 (Shows), 全部类型 (All types).
# 本地生活 (Local life) sheet is unchanged.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $RowValues) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $RowValues[$row]
    }
}

# Sheet: 展览
Set-FValues "展览" @{
    3  = 527
    5  = 501
    9  = 1007
    10 = 799
    11 = 231
    12 = 56
    15 = 271
    18 = 1323
    20 = 846
    21 = 1160
    22 = 2846
    24 = 688
    25 = 188
    28 = 1001
    29 = 347
    30 = 2999
    31 = 571
    32 = 531
}

# Sheet: 演出
Set-FValues "演出" @{
    3 = 519
    9 = 41
}

# Sheet: 全部类型
Set-FValues "全部类型" @{
    4  = 527
    7  = 501
    8  = 519
    16 = 1007
    17 = 799
    18 = 231
    20 = 56
    21 = 41
    27 = 271
    30 = 1323
    32 = 846
    33 = 1160
    34 = 2846
    36 = 688
    37 = 188
    42 = 1001
    43 = 347
    44 = 2999
    45 = 571
    46 = 531
}

Write-Host "F-column updates applied"
